$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Manufacturer 1" .. "Supplier Subtotal 1" columns (old F:L) in one
# shot, then remove the "Line #" column (old A) -- so Name/Description/
# Designator/Quantity shift left into A:D.
$ws.Range("F1:L1").EntireColumn.Delete() | Out-Null
$ws.Range("A1:A1").EntireColumn.Delete() | Out-Null

# Re-apply the column widths used for the re-exported layout.
$ws.Range("A:A").ColumnWidth = 15.6
$ws.Range("B:D").ColumnWidth = 18.8
